# TradingModel_v2 - 2021/11/19 data updated
# Row 12 (2021-11-18) was previously the last row and carried the
# "last row" date style (YYYY-MM-DD). Since a new row is being appended,
# row 12 reverts to the regular data-row date style (YYYY-MM-DD HH:MM:SS)
# and the new row 13 (2021-11-19) takes over the "last row" style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: revert from "last row" style back to the normal data-row style.
$ws.Range("A12").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 13: new last row with the "last row" date style.
$ws.Range("A13").NumberFormat = "YYYY-MM-DD"
$ws.Range("A13").Value = 44519
$ws.Range("B13").Value = 61854.9
